$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 1 (headers): shift the old "ExpTitle" header from G1 out to O1, and
# fill G1:Q1 with the new set of headers used by the two new test cases.
# First stamp the bold/filled header style across the whole new range (Q1
# ends up with the style but no value, same as the rest of the header row).
$ws.Range("A1").Copy($ws.Range("G1:Q1")) | Out-Null
$ws.Range("G1:Q1").ClearContents() | Out-Null

$ws.Range("G1").Value = "Lead Source"
$ws.Range("H1").Value = "No of employees"
$ws.Range("I1").Value = "Mobile"
$ws.Range("J1").Value = "Lead Status"
$ws.Range("K1").Value = "City"
$ws.Range("L1").Value = "EmailWrong"
$ws.Range("M1").Value = "EmailRight"
$ws.Range("N1").Value = "Assigned to"
$ws.Range("O1").Value = "ExpTitle"
$ws.Range("P1").Value = "ExpAlert1"

# --- Row 2: the old "vtiger CRM..." expected title now lives under O2.
$ws.Range("O2").Value = "vtiger CRM - Commercial Open Source CRM"
$ws.Range("G2").ClearContents() | Out-Null

# --- Row 7: new test case TC06_CreateLeadWithMandatoryFields
$ws.Range("A7").Value = "TC06_CreateLeadWithMandatoryFields"
$ws.Range("B7").Value = "admin"
$ws.Range("C7").Value = "admin"
$ws.Range("E7").Value = "Mishra"
$ws.Range("F7").Value = "Google"
$ws.Range("P7").Value = "Last Name cannot be empty"
$ws.Range("P7").Font.Name = "Courier New"
$ws.Range("P7").Font.Size = 10
$ws.Range("P7").Font.Color = 0
$ws.Range("P7").VerticalAlignment = -4108

# --- Row 8: new test case TC07_EditLead_DeleteLead
$ws.Range("A8").Value = "TC07_EditLead_DeleteLead"
$ws.Range("B8").Value = "admin"
$ws.Range("C8").Value = "admin"
$ws.Range("D8").Value = "Devendra"
$ws.Range("E8").Value = "Mishra"
$ws.Range("F8").Value = "Google"
$ws.Range("G8").Value = "Partner"
$ws.Range("H8").Value = 40
$ws.Range("I8").Value = 776127892
$ws.Range("J8").Value = 4
$ws.Range("K8").Value = "Pune"
$ws.Range("L8").Value = "hgdhjw"
$ws.Range("M8").Value = "ddhg@gmail.com"
$ws.Hyperlinks.Add($ws.Range("M8"), "mailto:ddhg@gmail.com") | Out-Null
$ws.Range("N8").Value = "Team"
$ws.Range("P8").Value = "Please enter a valid Email"

# --- Column widths (best-effort match of the author's bestFit autosize)
$ws.Columns.Item(1).ColumnWidth = 32.666666666666664
$ws.Columns.Item(9).ColumnWidth = 9.166666666666666
$ws.Columns.Item(12).ColumnWidth = 10.333333333333332
$ws.Columns.Item(13).ColumnWidth = 14.333333333333332
$ws.Columns.Item(14).ColumnWidth = 9.833333333333332
$ws.Columns.Item(15).ColumnWidth = 36.166666666666664
$ws.Columns.Item(16).ColumnWidth = 23.333333333333336
$ws.Columns.Item(17).ColumnWidth = 22.333333333333336

# --- Final selection lands on M1, matching the saved cursor position.
$ws.Range("M1").Select() | Out-Null
